# Auto-generated edit script: update crypto Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.590.97'
$ws.Range("E2").Value = '  +1.01%  '
$ws.Range("D3").Value = '3.395.45'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'576.04"
$ws.Range("E5").Value = '  +0.83%  '
$ws.Range("D6").Value = "'142.00"
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = "'0.473"
$ws.Range("E8").Value = '  -0.45%  '
$ws.Range("D9").Value = "'7.65"
$ws.Range("E9").Value = '  +0.61%  '
$ws.Range("E10").Value = '  -0.77%  '
$ws.Range("E11").Value = '  -1.62%  '
$ws.Range("D12").Value = '3.975.86'
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").Value = "'27.99"
$ws.Range("E14").Value = '  +1.29%  '
$ws.Range("D15").Value = '3.392.85'
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("E16").Value = '  -0.27%  '
$ws.Range("D17").Value = '61.651.54'
$ws.Range("E17").Value = '  +0.95%  '
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("E19").Value = '  -0.77%  '
$ws.Range("D20").Value = "'9.05"
$ws.Range("E20").Value = '  +1.46%  '
$ws.Range("D21").Value = "'387.03"
$ws.Range("E21").Value = '  +1.50%  '
$ws.Range("D22").Value = "'74.67"
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("E23").Value = '  -1.05%  '
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("E25").Value = '  -3.82%  '
$ws.Range("D26").Value = "'0.195"
$ws.Range("E26").Value = '  +7.85%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("D28").Value = "'7.41"
$ws.Range("E28").Value = '  +1.15%  '
$ws.Range("D29").Value = "'7.98"
$ws.Range("E29").Value = '  +0.35%  '
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("D31").Value = "'1.41"
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").Value = "'23.33"
$ws.Range("E33").Value = '  -0.17%  '
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("D35").Value = "'168.35"
$ws.Range("E35").Value = '  +1.16%  '
$ws.Range("E36").Value = '  +1.00%  '
$ws.Range("D37").Value = '3.431.00'
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("E38").Value = '  +0.31%  '
$ws.Range("D39").Value = "'0.0764"
$ws.Range("E39").Value = '  -0.43%  '
$ws.Range("D40").Value = "'26.44"
$ws.Range("E40").Value = '  -3.58%  '
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("E42").Value = '  +0.76%  '
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("E44").Value = '  +1.87%  '
$ws.Range("D45").Value = '2.457.40'
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").Value = "'22.69"
$ws.Range("E46").Value = '  -1.08%  '
$ws.Range("E47").Value = '  -1.35%  '
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("E49").Value = '  -1.58%  '
$ws.Range("D50").Value = "'2.03"
$ws.Range("E50").Value = '  -4.29%  '
$ws.Range("E51").Value = '  -1.03%  '
